# Insert a new price-record row at row 88 (pushes existing rows 88..207
# down to 89..208) and populate it with a new "Albahaca" observation,
# matching the weekly Fruta/Hortaliza data refresh described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 88..207 down to 89..208, leaving a blank row 88 behind.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new observation.
$ws.Cells.Item(88, 1).Value = 3
$ws.Cells.Item(88, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(88, 3).Value = 'Coquimbo'
$ws.Cells.Item(88, 4).Value = 44930
$ws.Cells.Item(88, 5).Value = 5
$ws.Cells.Item(88, 6).Value = 100112052
$ws.Cells.Item(88, 7).Value = 'Albahaca'
$ws.Cells.Item(88, 8).Value = 'Sin especificar'
$ws.Cells.Item(88, 9).Value = 'Primera'
$ws.Cells.Item(88, 10).Value = 115
$ws.Cells.Item(88, 11).Value = 5500
$ws.Cells.Item(88, 12).Value = 6000
$ws.Cells.Item(88, 13).Value = 5739
$ws.Cells.Item(88, 14).Value = '$/docena de matas'
$ws.Cells.Item(88, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(88, 16).Value = 956
$ws.Cells.Item(88, 17).Value = 6
$ws.Cells.Item(88, 18).Value = 'Hortaliza'
